# "Generate Report for Archive"
# - Status text "Ready for handoff" -> "In Translation" on all 3 sheets
#   (Overview!E2:F2, zh-cn!C2, de-de!C2 all point at the same shared
#   "Status" string, so every occurrence has to be updated).
# - Narrow the "Status" column (Overview cols E & F, zh-cn/de-de col C)
#   from ~17.216 chars to ~13.410 chars.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update the Status value everywhere it appears ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status columns ---
# Target stored column width is ~13.4101845877511 characters.
$newStatusWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
